$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store plain text that looks like a
# number or a signed percentage (e.g. "542.71", "  +0.05%  "). Excel's
# normal typed-entry parsing would silently turn these into floating point
# Numbers, so each such cell is marked Text (NumberFormat "@") before its
# new literal string is written, to match the source data exactly.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.108.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.347.85'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.59'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.104'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.11%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.81'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.765.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.035.12'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.361.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.30'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '329.09'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '63.03'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.165'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.33'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.33'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -7.63%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.20'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0738'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.36'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.17'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.61'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.15'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '141.25'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.378'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '289.62'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.64'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0948'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0511'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.99'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.381'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('B51').Value = 'ZEEBU'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.68'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.27%  '
